$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.862.64"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "1.751.17"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "236.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9987"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5132"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.60%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.40"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2669"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06179"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.88%  "
$ws.Range("D11").Value = "1.782.58"
$ws.Range("E11").Value = "  +2.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06944"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.53%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.39"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.21%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6282"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +12.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.479"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "77.97"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.95%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9996"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.15%  "
$ws.Range("D19").Value = "25.869.78"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.61"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006653"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "1.985.46"
$ws.Range("E22").Value = "  +1.19%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.058"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.71%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.264"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "136.30"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.481"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("E28").Value = "  +3.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.778"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.84%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08278"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.95%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.688"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.402"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04385"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.637"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9986"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.75%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6051"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.60%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.665"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01557"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.59%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.933"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.12%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9989"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.16"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.89%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3838"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7477"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.22%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.887"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05488"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.28%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1099"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.42%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.961"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.03%  "
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "52.62"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.004"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.22%  "
